{"js": "// Apply the four content edits described by the diff:\n//  1. Update the cached DATE field result from 10/13/2025 to 10/23/2025.\n//  2. Change \"adjustable,\" to \"adjustable;\" in the Scope/limitations bullet.\n//  3. Replace the Vina/GNINA explanation sentence with the new wording.\n//  4. Drop \"AutoDock Vina, \" from the Docking/Screening tool list.\n\nconst body = context.document.body;\n\n// 1) Date field cached value: 10/13/2025 -> 10/23/2025\nconst dateHits = body.search(\"10/13/2025\", { matchCase: true });\ndateHits.load(\"text\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"10/23/2025\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"adjustable,\" -> \"adjustable;\"\nconst commaHits = body.search(\"grid box size adjustable,\", { matchCase: true });\ncommaHits.load(\"text\");\nawait context.sync();\nif (commaHits.items.length > 0) {\n  commaHits.items[0].insertText(\"grid box size adjustable;\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Replace the Vina/GNINA sentence\nconst oldSentence =\n  \"Vina can be used for quick and basic docking, while GNINA can be used for accurate pose prediction and binding affinity estimation (we can choose either, Gnina = Vina + machine learning)\";\nconst newSentence =\n  \"since GNINA is based on Autodock Vina and smina and has a machine learning feature, this will be the primary docking tool in this workflow.\";\nconst sentenceHits = body.search(oldSentence, { matchCase: true });\nsentenceHits.load(\"text\");\nawait context.sync();\nif (sentenceHits.items.length > 0) {\n  sentenceHits.items[0].insertText(newSentence, \"Replace\");\n  await context.sync();\n}\n\n// 4) Docking/Screening tool list: drop \"AutoDock Vina, \"\nconst toolHits = body.search(\"AutoDock Vina, Gnina\", { matchCase: true });\ntoolHits.load(\"text\");\nawait context.sync();\nif (toolHits.items.length > 0) {\n  toolHits.items[0].insertText(\"Gnina\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the four content edits described by the diff:\n#  1. Update the cached DATE field result from 10/13/2025 to 10/23/2025.\n#  2. Change \"adjustable,\" to \"adjustable;\" in the Scope/limitations bullet.\n#  3. Replace the Vina/GNINA explanation sentence with the new wording.\n#  4. Drop \"AutoDock Vina, \" from the Docking/Screening tool list.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param(\n        [string]$OldText,\n        [string]$NewText\n    )\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $OldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $NewText\n    $wdFindContinue = 1\n    $wdReplaceAll = 2\n    return $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n\n# 1) Date field cached value: 10/13/2025 -> 10/23/2025\nReplace-DocText \"10/13/2025\" \"10/23/2025\" | Out-Null\n\n# 2) \"adjustable,\" -> \"adjustable;\"\nReplace-DocText \"grid box size adjustable,\" \"grid box size adjustable;\" | Out-Null\n\n# 3) Replace the Vina/GNINA sentence\n$oldSentence = \"Vina can be used for quick and basic docking, while GNINA can be used for accurate pose prediction and binding affinity estimation (we can choose either, Gnina = Vina + machine learning)\"\n$newSentence = \"since GNINA is based on Autodock Vina and smina and has a machine learning feature, this will be the primary docking tool in this workflow.\"\nReplace-DocText $oldSentence $newSentence | Out-Null\n\n# 4) Docking/Screening tool list: drop \"AutoDock Vina, \"\nReplace-DocText \"AutoDock Vina, Gnina\" \"Gnina\" | Out-Null\n"}
